$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.429.70"
$ws.Range("E2").Value = "  -1.68%  "
$ws.Range("D3").Value = "2.574.34"
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.69%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.43%  "
$ws.Range("D9").Value = "2.581.56"
$ws.Range("E9").Value = "  -2.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.163"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.96%  "
$ws.Range("E13").Value = "  +2.29%  "
$ws.Range("D14").Value = "3.024.19"
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("D15").Value = "59.411.42"
$ws.Range("E15").Value = "  -1.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000137"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "2.573.44"
$ws.Range("E18").Value = "  -2.55%  "
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.482"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("E27").Value = "  -2.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("D29").Value = "0.0₃0771"
$ws.Range("E29").Value = "  -2.61%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("E32").Value = "  -2.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +2.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.909"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.854"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.81%  "
$ws.Range("E40").Value = "  -1.99%  "
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "287.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "135.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.55%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0968"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("E46").Value = "  -1.65%  "
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("E48").Value = "  -2.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0234"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").Value = "1.950.17"
$ws.Range("E51").Value = "  -0.62%  "
